$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MappingTable")

# --- Step 1: fill in the previously-empty B/C/D mapping columns for the
# Texit*/Twall*/Tpack* rows (6,7,8) and update their MappedVariable (E) names.
$ws.Range("B6").Value = "exit_gas_temperatures"
$ws.Range("C6").Value = "exit_gas_temperature_units"
$ws.Range("D6").Value = "filter"
$ws.Range("E6").Value = "EGT"

$ws.Range("B7").Value = "wall_left_temperatures"
$ws.Range("C7").Value = "wall_temperature_units"
$ws.Range("D7").Value = "chamber"
$ws.Range("E7").Value = "TWALL"

$ws.Range("B8").Value = "filter_temperatures"
$ws.Range("C8").Value = "filter_temperature_units"
$ws.Range("D8").Value = "filter"
$ws.Range("E8").Value = "TPACK"

# --- Step 2: fill in the previously-empty B/C/D mapping columns for the
# Qpack* (row 25) and Epack* (row 27) rows; their MappedVariable (E) stays.
$ws.Range("B25").Value = "filter_absorption_rates"
$ws.Range("C25").Value = "filter_absorption_rate_units"
$ws.Range("D25").Value = "filter"

$ws.Range("B27").Value = "filter_energies"
$ws.Range("C27").Value = "filter_energy_units"
$ws.Range("D27").Value = "filter"

# --- Step 3: delete the rows that are no longer part of the mapping table.
# Deleting from the bottom up keeps the row numbers above stable.
$ws.Rows("37:37").Delete()
$ws.Rows("29:29").Delete()
$ws.Rows("10:24").Delete()
